$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "L1cam"
$ws.Range("C2").Value = "Erbb2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 18.04537966666667
$ws.Range("H2").Value = 54.13613900000001
$ws.Range("I2").Value = 0.6797959733292525
$ws.Range("J2").Value = 0.6797959733292525
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 1.720171333333333
$ws.Range("N2").Value = 5.160514
$ws.Range("O2").Value = 0.1961456356393658
$ws.Range("P2").Value = 0.1961456356393658
$ws.Range("Q2").Value = 31.04114480171623
$ws.Range("R2").Value = 279.370303215446
$ws.Range("S2").Value = 0.1333390132937476
$ws.Range("T2").Value = 0.1333390132937476

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "L1cam"
$ws.Range("C3").Value = "Erbb2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 18.04537966666667
$ws.Range("H3").Value = 54.13613900000001
$ws.Range("I3").Value = 0.6797959733292525
$ws.Range("J3").Value = 0.6797959733292525
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 4.076388666666666
$ws.Range("N3").Value = 12.229166
$ws.Range("O3").Value = 0.464817562438416
$ws.Range("P3").Value = 0.464817562438416
$ws.Range("Q3").Value = 73.55998115889712
$ws.Range("R3").Value = 662.0398304300741
$ws.Range("S3").Value = 0.3159811072783536
$ws.Range("T3").Value = 0.3159811072783535

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "L1cam"
$ws.Range("C4").Value = "Erbb2"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 18.04537966666667
$ws.Range("H4").Value = 54.13613900000001
$ws.Range("I4").Value = 0.6797959733292525
$ws.Range("J4").Value = 0.6797959733292525
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 2.973308
$ws.Range("N4").Value = 8.919924
$ws.Range("O4").Value = 0.3390368019222182
$ws.Range("P4").Value = 0.3390368019222182
$ws.Range("Q4").Value = 53.65447172593733
$ws.Range("R4").Value = 482.8902455334361
$ws.Range("S4").Value = 0.2304758527571513
$ws.Range("T4").Value = 0.2304758527571513

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "L1cam"
$ws.Range("C5").Value = "Erbb2"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.6001993333333334
$ws.Range("H5").Value = 1.800598
$ws.Range("I5").Value = 0.02261039099934159
$ws.Range("J5").Value = 0.02261039099934159
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 1.720171333333333
$ws.Range("N5").Value = 5.160514
$ws.Range("O5").Value = 0.1961456356393658
$ws.Range("P5").Value = 0.1961456356393658
$ws.Range("Q5").Value = 1.032445687485778
$ws.Range("R5").Value = 9.292011187372001
$ws.Range("S5").Value = 0.004434929514620451
$ws.Range("T5").Value = 0.004434929514620451

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "L1cam"
$ws.Range("C6").Value = "Erbb2"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.6001993333333334
$ws.Range("H6").Value = 1.800598
$ws.Range("I6").Value = 0.02261039099934159
$ws.Range("J6").Value = 0.02261039099934159
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 4.076388666666666
$ws.Range("N6").Value = 12.229166
$ws.Range("O6").Value = 0.464817562438416
$ws.Range("P6").Value = 0.464817562438416
$ws.Range("Q6").Value = 2.446645760140889
$ws.Range("R6").Value = 22.019811841268
$ws.Range("S6").Value = 0.01050970683009346
$ws.Range("T6").Value = 0.01050970683009346

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "L1cam"
$ws.Range("C7").Value = "Erbb2"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.6001993333333334
$ws.Range("H7").Value = 1.800598
$ws.Range("I7").Value = 0.02261039099934159
$ws.Range("J7").Value = 0.02261039099934159
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 2.973308
$ws.Range("N7").Value = 8.919924
$ws.Range("O7").Value = 0.3390368019222182
$ws.Range("P7").Value = 0.3390368019222182
$ws.Range("Q7").Value = 1.784577479394667
$ws.Range("R7").Value = 16.061197314552
$ws.Range("S7").Value = 0.007665754654627682
$ws.Range("T7").Value = 0.007665754654627681

# Row 8
$ws.Range("A8").Value = "M2"
$ws.Range("B8").Value = "L1cam"
$ws.Range("C8").Value = "Erbb2"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 4.690054333333333
$ws.Range("H8").Value = 14.070163
$ws.Range("I8").Value = 0.1766812397072912
$ws.Range("J8").Value = 0.1766812397072912
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 1.720171333333333
$ws.Range("N8").Value = 5.160514
$ws.Range("O8").Value = 0.1961456356393658
$ws.Range("P8").Value = 0.1961456356393658
$ws.Range("Q8").Value = 8.067697015975778
$ws.Range("R8").Value = 72.609273143782
$ws.Range("S8").Value = 0.03465525406793778
$ws.Range("T8").Value = 0.03465525406793778

# Row 9
$ws.Range("A9").Value = "M2"
$ws.Range("B9").Value = "L1cam"
$ws.Range("C9").Value = "Erbb2"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 4.690054333333333
$ws.Range("H9").Value = 14.070163
$ws.Range("I9").Value = 0.1766812397072912
$ws.Range("J9").Value = 0.1766812397072912
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 4.076388666666666
$ws.Range("N9").Value = 12.229166
$ws.Range("O9").Value = 0.464817562438416
$ws.Range("P9").Value = 0.464817562438416
$ws.Range("Q9").Value = 19.11848433045089
$ws.Range("R9").Value = 172.066358974058
$ws.Range("S9").Value = 0.08212454316934056
$ws.Range("T9").Value = 0.08212454316934056

# Row 10
$ws.Range("A10").Value = "M2"
$ws.Range("B10").Value = "L1cam"
$ws.Range("C10").Value = "Erbb2"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 4.690054333333333
$ws.Range("H10").Value = 14.070163
$ws.Range("I10").Value = 0.1766812397072912
$ws.Range("J10").Value = 0.1766812397072912
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 2.973308
$ws.Range("N10").Value = 8.919924
$ws.Range("O10").Value = 0.3390368019222182
$ws.Range("P10").Value = 0.3390368019222182
$ws.Range("Q10").Value = 13.94497606973467
$ws.Range("R10").Value = 125.504784627612
$ws.Range("S10").Value = 0.05990144247001283
$ws.Range("T10").Value = 0.05990144247001283

# Row 11
$ws.Range("A11").Value = "sCs"
$ws.Range("B11").Value = "L1cam"
$ws.Range("C11").Value = "Erbb2"
$ws.Range("D11").Value = "ECs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 3.209654333333333
$ws.Range("H11").Value = 9.628962999999999
$ws.Range("I11").Value = 0.1209123959641148
$ws.Range("J11").Value = 0.1209123959641148
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 1.720171333333333
$ws.Range("N11").Value = 5.160514
$ws.Range("O11").Value = 0.1961456356393658
$ws.Range("P11").Value = 0.1961456356393658
$ws.Range("Q11").Value = 5.52115537410911
$ws.Range("R11").Value = 49.690398366982
$ws.Range("S11").Value = 0.02371643876305998
$ws.Range("T11").Value = 0.02371643876305998

# Row 12
$ws.Range("A12").Value = "sCs"
$ws.Range("B12").Value = "L1cam"
$ws.Range("C12").Value = "Erbb2"
$ws.Range("D12").Value = "FAPs"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 3.209654333333333
$ws.Range("H12").Value = 9.628962999999999
$ws.Range("I12").Value = 0.1209123959641148
$ws.Range("J12").Value = 0.1209123959641148
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 4.076388666666666
$ws.Range("N12").Value = 12.229166
$ws.Range("O12").Value = 0.464817562438416
$ws.Range("P12").Value = 0.464817562438416
$ws.Range("Q12").Value = 13.08379854831755
$ws.Range("R12").Value = 117.754186934858
$ws.Range("S12").Value = 0.05620220516062841
$ws.Range("T12").Value = 0.0562022051606284

# Row 13
$ws.Range("A13").Value = "sCs"
$ws.Range("B13").Value = "L1cam"
$ws.Range("C13").Value = "Erbb2"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 3.209654333333333
$ws.Range("H13").Value = 9.628962999999999
$ws.Range("I13").Value = 0.1209123959641148
$ws.Range("J13").Value = 0.1209123959641148
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 2.973308
$ws.Range("N13").Value = 8.919924
$ws.Range("O13").Value = 0.3390368019222182
$ws.Range("P13").Value = 0.3390368019222182
$ws.Range("Q13").Value = 9.543290906534665
$ws.Range("R13").Value = 85.88961815881198
$ws.Range("S13").Value = 0.04099375204042641
$ws.Range("T13").Value = 0.0409937520404264
